$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7-76 down to 8-77.
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new weekly record (same fixed attributes as the
# surrounding rows, with updated date / volume / price / origin data).
$ws.Range("A7").Value = 11
$ws.Range("B7").Value = "Vega Monumental Concepción"
$ws.Range("C7").Value = "Bíobío"
$ws.Range("D7").Value = 44635
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 100112001
$ws.Range("G7").Value = "Berenjena"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 140
$ws.Range("K7").Value = 7500
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 7786
$ws.Range("N7").Value = "$/caja 60 unidades"
$ws.Range("O7").Value = "Región de Arica y Parinacota"
$ws.Range("P7").Value = 130
$ws.Range("Q7").Value = 60
$ws.Range("R7").Value = "Hortaliza"

# Make sure the date cell keeps the date-formatted style used throughout column D.
$ws.Range("D7").NumberFormat = $ws.Range("D8").NumberFormat
